$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 11
$ws.Cells.Item($row, 1).Value = "Globo"
$ws.Cells.Item($row, 2).Value = "RJ TV 1"
$ws.Cells.Item($row, 3).Value = "Esportes"
$ws.Cells.Item($row, 4).Value = "2025-04-01T11:43"
$ws.Cells.Item($row, 5).Value = "Neutro"
$ws.Cells.Item($row, 6).Value = "Presidente destituído. Após assembleia, Laila Póvoa assume a presidência do Americano. "
